$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Isokov Eldor Fayzullayevich"
$ws.Range("B5").Value = "Mehnat muhofazasi va texnika xavfsizligi"
$ws.Range("C5").Value = "O'zbek tili"
$ws.Range("D5").Value = "Kunduzgi"
$ws.Range("E5").Value = "AD1460068"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "3140687256008"
$ws.Range("G5").Value = "Toshkent shahri"
$ws.Range("H5").Value = "Mirzo Ulugʻbek tumani"
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "998971300087"
